$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A130").Value = 'What’s the track limit for an ODF?'
$ws.Range("B130").Value = 'llama3.2:latest'
$ws.Range("C130").Value = 'The maximum number of tracks that can be specified within a single ODF file is 200 tracks at most. It''s recommended to use OIF files instead, which allows you to save only the specific depth range of interest, making data management more efficient.'

$ws.Range("A131").Value = 'How to open the GEO application on Windows?'
$ws.Range("B131").Value = 'llama3.2:latest'
$ws.Range("C131").Value = 'To open the GEO application on Windows, follow these steps:
1. Click the Windows button and select Geologix GEO Suite 8, then GEO.
2. In Windows Explorer, locate the installation location and double-click the Geo.exe file.
Alternatively, you can also open the application by clicking the program icon (GEO for Windows) in the Windows menu.'

$ws.Range("A132").Value = 'How do you define new curve pens?'
$ws.Range("B132").Value = 'llama3.2:latest'
$ws.Range("C132").Value = 'To define a new curve pen, follow these steps:
1. On the Curve Settings tab, select Curve Settings.
2. The Curve Settings sidebar will open.
3. On the Display tab, open the Pens node.
4. Select the "Use Customized Pens" double arrow.
5. The Pen Definition dialog box will open.
6. Select an unused number in the "Select Pen" box.
7. Enter a Pen Name (e.g., DotDashDash).
8. Select the Pen Type for both the display and printing pens.
Note: A pen definition consists of up to five segments of dot-combinations, with up to ten pixels on and up to ten pixels off per segment.'
